# "change pca to svd"
#
# Summary of the edit (per the authoritative OOXML diff):
#  1. Sheet1!G7 rich text: fold the old 2nd run into the (unstyled) 1st run and
#     append a brand-new, regular-weight 2nd run describing the new SVD results.
#  2. Sheet1!K7: turn the plain string into two runs - the original text plus a
#     comma, and a new regular-weight run describing the new NMF result.
#  3. Sheet1 becomes the active sheet/tab (was Sheet2); selection moves to G6;
#     column G gets wider.
#  4. Sheet2 stops being the active tab (selection itself is untouched).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. G7: pca -> pca + new svd run ------------------------------------
$g7 = "pca_k=4:8.6%;kmeans_k=4:9.7%; pca_k=100:16.1%  pca_k=4(ordered):5.3%;pca_k=100(ordered):5.9% svd_k=4(ordered):8.6%,svd_k=100(ordered):8.6%"
$ws1.Range("G7").Value = $g7

$g7Run2Start = 94
$g7Run2Len = 45
$g7Run2 = $ws1.Range("G7").Characters($g7Run2Start, $g7Run2Len)
$g7Run2.Font.Name = "宋体"
$g7Run2.Font.Size = 11
$g7Run2.Font.Color = 0

# --- 2. K7: pca_k=100:15.5% -> + new nmf run ----------------------------
$k7 = "pca_k=100:15.5%,nmf_k=4(ordered):8.7%"
$ws1.Range("K7").Value = $k7

$k7Run2Start = 17
$k7Run2Len = 21
$k7Run2 = $ws1.Range("K7").Characters($k7Run2Start, $k7Run2Len)
$k7Run2.Font.Name = "宋体"
$k7Run2.Font.Size = 11
$k7Run2.Font.Color = 0

# --- 3. Column G gets wider ----------------------------------------------
$ws1.Columns.Item(7).ColumnWidth = 48.7

# --- 4. View/selection changes -------------------------------------------
# Sheet1 becomes the active sheet/tab, selection on G6 (was Sheet2 active,
# Sheet1 selection was G7).
$ws1.Activate()
$ws1.Range("G6").Select()

# Sheet2's selection (F5) is left untouched; it simply stops being the
# active tab once Sheet1 is activated above.
